$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 161877
$ws.Range("C4").Value = 152885
$ws.Range("C7").Value = 5.55
$ws.Range("C8").Value = 64.59
